# Add I0 and IF columns (I and J) mirroring existing IP (column H) data.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Headers
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)  # xlPasteFormats

# Data rows 2-38: I = 1, J = same value as H (IP)
for ($r = 2; $r -le 38; $r++) {
    $ipValue = $ws.Cells.Item($r, 8).Value2
    $ws.Cells.Item($r, 9).Value = 1
    $ws.Cells.Item($r, 10).Value = $ipValue
}
